# Update "want to go" counts (column F, header "想去人数") across the
# "展览" (rId1/sheet1), "演出" (rId2/sheet2) and "全部类型" (rId4/sheet4)
# sheets, matching the regenerated output for commit 456a3b4.
# ("本地生活" / sheet3 has no changes in this commit.)

$wb = $excel.ActiveWorkbook

# 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 113
$ws.Range("F4").Value = 108
$ws.Range("F5").Value = 291
$ws.Range("F6").Value = 178
$ws.Range("F7").Value = 1154
$ws.Range("F8").Value = 402
$ws.Range("F9").Value = 91
$ws.Range("F10").Value = 107
$ws.Range("F15").Value = 139
$ws.Range("F16").Value = 1322
$ws.Range("F18").Value = 188
$ws.Range("F21").Value = 672
$ws.Range("F22").Value = 1076
$ws.Range("F24").Value = 1936
$ws.Range("F25").Value = 2495
$ws.Range("F26").Value = 1260
$ws.Range("F28").Value = 189
$ws.Range("F29").Value = 358
$ws.Range("F30").Value = 797
$ws.Range("F32").Value = 964
$ws.Range("F33").Value = 118
$ws.Range("F35").Value = 746
$ws.Range("F36").Value = 368
$ws.Range("F37").Value = 576
$ws.Range("F38").Value = 715
$ws.Range("F39").Value = 313
$ws.Range("F40").Value = 208

# 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F14").Value = 366

# 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 113
$ws.Range("F6").Value = 108
$ws.Range("F7").Value = 291
$ws.Range("F8").Value = 178
$ws.Range("F11").Value = 1154
$ws.Range("F12").Value = 402
$ws.Range("F13").Value = 91
$ws.Range("F14").Value = 107
$ws.Range("F20").Value = 139
$ws.Range("F21").Value = 1322
$ws.Range("F23").Value = 188
$ws.Range("F26").Value = 1076
$ws.Range("F27").Value = 2495
$ws.Range("F29").Value = 1260
$ws.Range("F34").Value = 189
$ws.Range("F35").Value = 358
$ws.Range("F36").Value = 797
$ws.Range("F40").Value = 964
$ws.Range("F41").Value = 746
$ws.Range("F42").Value = 368
$ws.Range("F43").Value = 576
$ws.Range("F44").Value = 715
$ws.Range("F45").Value = 313
$ws.Range("F48").Value = 208
